$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,13

$arr[0,0] = 1.02
$arr[0,1] = 1.02144895620723
$arr[0,2] = 1.032082841100385
$arr[0,3] = 1.022325574799219
$arr[0,4] = 1.043107793232937
$arr[0,5] = 1
$arr[0,7] = 1.032046144391726
$arr[0,8] = 1.026640362567578
$arr[0,9] = 1.034889096898955
$arr[0,10] = 1.025160287473952
$arr[0,11] = 1.045882590493306
$arr[0,12] = 1.012943769618517
$arr[1,0] = 1.02
$arr[1,1] = 1.022270124313473
$arr[1,2] = 1.032696184933035
$arr[1,3] = 1.023019017595175
$arr[1,4] = 1.043883397499351
$arr[1,5] = 1
$arr[1,7] = 1.032156974800357
$arr[1,8] = 1.027099822265489
$arr[1,9] = 1.03531166212589
$arr[1,10] = 1.025660681850534
$arr[1,11] = 1.04646925513283
$arr[1,12] = 1.013097350888327
$arr[2,0] = 1.02
$arr[2,1] = 1.022802240274017
$arr[2,2] = 1.033093544063434
$arr[2,3] = 1.023468775297714
$arr[2,4] = 1.044386089551724
$arr[2,5] = 1
$arr[2,7] = 1.032227555124078
$arr[2,8] = 1.027397242483795
$arr[2,9] = 1.035584881179608
$arr[2,10] = 1.025984848788732
$arr[2,11] = 1.046849030052153
$arr[2,12] = 1.013196723134852
$arr[3,0] = 1.02
$arr[3,1] = 1.023026123285496
$arr[3,2] = 1.033260707940521
$arr[3,3] = 1.023658103816721
$arr[3,4] = 1.044597616450248
$arr[3,5] = 1
$arr[3,7] = 1.032256954807405
$arr[3,8] = 1.02752230511545
$arr[3,9] = 1.03569969103567
$arr[3,10] = 1.026121217695064
$arr[3,11] = 1.047008724808799
$arr[3,12] = 1.013238497515702
$arr[4,0] = 1.02
$arr[4,1] = 1.023063724818456
$arr[4,2] = 1.033288782125682
$arr[4,3] = 1.023689907534974
$arr[4,4] = 1.044633144128626
$arr[4,5] = 1
$arr[4,7] = 1.032261875155881
$arr[4,8] = 1.027543305236186
$arr[4,9] = 1.035718965058201
$arr[4,10] = 1.026144119818839
$arr[4,11] = 1.047035540400724
$arr[4,12] = 1.013245511503451
$arr[5,0] = 1.02
$arr[5,1] = 1.022805231098943
$arr[5,2] = 1.033095777269444
$arr[5,3] = 1.023471304133233
$arr[5,4] = 1.044388915220373
$arr[5,5] = 1
$arr[5,7] = 1.032227949034873
$arr[5,8] = 1.02739891347171
$arr[5,9] = 1.035586415478413
$arr[5,10] = 1.025986670608348
$arr[5,11] = 1.046851163756695
$arr[5,12] = 1.013197281333514
$arr[6,0] = 1.02
$arr[6,1] = 1.021726314303236
$arr[6,2] = 1.03229002200949
$arr[6,3] = 1.022559707684588
$arr[6,4] = 1.043369740364203
$arr[6,5] = 1
$arr[6,7] = 1.032083834493352
$arr[6,8] = 1.026795613667417
$arr[6,9] = 1.035031947196034
$arr[6,10] = 1.025329318895486
$arr[6,11] = 1.04608082158299
$arr[6,12] = 1.012995673890493
$arr[7,0] = 1.02
$arr[7,1] = 1.019831063706675
$arr[7,2] = 1.030873975439597
$arr[7,3] = 1.02096151396548
$arr[7,4] = 1.041580232261256
$arr[7,5] = 1
$arr[7,7] = 1.031821230106229
$arr[7,8] = 1.025733493250462
$arr[7,9] = 1.034053366144892
$arr[7,10] = 1.024173942485544
$arr[7,11] = 1.044724713415577
$arr[7,12] = 1.012640398353874
$arr[8,0] = 1.02
$arr[8,1] = 1.018571660598828
$arr[8,2] = 1.029932612626286
$arr[8,3] = 1.019901649180114
$arr[8,4] = 1.040391664592958
$arr[8,5] = 1
$arr[8,7] = 1.031640383440009
$arr[8,8] = 1.025026145422703
$arr[8,9] = 1.033400030254269
$arr[8,10] = 1.023405768827599
$arr[8,11] = 1.043821645300917
$arr[8,12] = 1.012403565472651
$arr[9,0] = 1.02
$arr[9,1] = 1.018027319015214
$arr[9,2] = 1.029525649421245
$arr[9,3] = 1.019444067946622
$arr[9,4] = 1.03987808160762
$arr[9,5] = 1
$arr[9,7] = 1.031560713637606
$arr[9,8] = 1.024720046576956
$arr[9,9] = 1.033116921592311
$arr[9,10] = 1.023073652120154
$arr[9,11] = 1.043430867757603
$arr[9,12] = 1.012301024544981
$arr[10,0] = 1.02
$arr[10,1] = 1.01782527656936
$arr[10,2] = 1.029374585273956
$arr[10,3] = 1.019274306282414
$arr[10,4] = 1.039687477699881
$arr[10,5] = 1
$arr[10,7] = 1.031530916816233
$arr[10,8] = 1.024606377357979
$arr[10,9] = 1.033011732343983
$arr[10,10] = 1.022950366924231
$arr[10,11] = 1.043285756021852
$arr[10,12] = 1.012262938171469
$arr[11,0] = 1.02
$arr[11,1] = 1.017868608531903
$arr[11,2] = 1.029406984484453
$arr[11,3] = 1.019310711455344
$arr[11,4] = 1.039728355442533
$arr[11,5] = 1
$arr[11,7] = 1.031537317551394
$arr[11,8] = 1.024630758433545
$arr[11,9] = 1.033034297132993
$arr[11,10] = 1.022976808475726
$arr[11,11] = 1.043316881122013
$arr[11,12] = 1.012271107734215
$arr[12,0] = 1.02
$arr[12,1] = 1.018010615048159
$arr[12,2] = 1.029513160358083
$arr[12,3] = 1.019430031213799
$arr[12,4] = 1.039862322877259
$arr[12,5] = 1
$arr[12,7] = 1.031558254778604
$arr[12,8] = 1.024710650033755
$arr[12,9] = 1.033108227221588
$arr[12,10] = 1.023063459733028
$arr[12,11] = 1.043418871940563
$arr[12,12] = 1.012297876271969
$arr[13,0] = 1.02
$arr[13,1] = 1.018098129913134
$arr[13,2] = 1.029578592142556
$arr[13,3] = 1.019503575196387
$arr[13,4] = 1.039944886391018
$arr[13,5] = 1
$arr[13,7] = 1.031571127896752
$arr[13,8] = 1.024759877837993
$arr[13,9] = 1.033153774045872
$arr[13,10] = 1.023116858777681
$arr[13,11] = 1.043481717264279
$arr[13,12] = 1.012314369517279
$arr[14,0] = 1.02
$arr[14,1] = 1.018607807730485
$arr[14,2] = 1.029959635363639
$arr[14,3] = 1.019932045892555
$arr[14,4] = 1.040425772236047
$arr[14,5] = 1
$arr[14,7] = 1.031645642227681
$arr[14,8] = 1.025046464260219
$arr[14,9] = 1.033418814945071
$arr[14,10] = 1.023427821139843
$arr[14,11] = 1.043847585489829
$arr[14,12] = 1.012410371014331
$arr[15,0] = 1.02
$arr[15,1] = 1.018927781002918
$arr[15,2] = 1.030198829966735
$arr[15,3] = 1.020201176323311
$arr[15,4] = 1.0407277085098
$arr[15,5] = 1
$arr[15,7] = 1.031692018957285
$arr[15,8] = 1.025226283493761
$arr[15,9] = 1.033585012924546
$arr[15,10] = 1.023623016548
$arr[15,11] = 1.044077155046919
$arr[15,12] = 1.012470593052525
$arr[16,0] = 1.02
$arr[16,1] = 1.019114511233793
$arr[16,2] = 1.030338411014094
$arr[16,3] = 1.020358285404352
$arr[16,4] = 1.040903926369774
$arr[16,5] = 1
$arr[16,7] = 1.031718938257475
$arr[16,8] = 1.025331187003508
$arr[16,9] = 1.033681933043638
$arr[16,10] = 1.023736919599572
$arr[16,11] = 1.044211083749947
$arr[16,12] = 1.012505720409749
$arr[17,0] = 1.02
$arr[17,1] = 1.019178197516361
$arr[17,2] = 1.030386015149607
$arr[17,3] = 1.020411877506268
$arr[17,4] = 1.04096402960058
$arr[17,5] = 1
$arr[17,7] = 1.03172809472325
$arr[17,8] = 1.025366959399883
$arr[17,9] = 1.033714976803658
$arr[17,10] = 1.023775765841136
$arr[17,11] = 1.044256754103992
$arr[17,12] = 1.012517698061661
$arr[18,0] = 1.02
$arr[18,1] = 1.018893441035598
$arr[18,2] = 1.030173160131013
$arr[18,3] = 1.020172287729444
$arr[18,4] = 1.040695302879461
$arr[18,5] = 1
$arr[18,7] = 1.031687056765628
$arr[18,8] = 1.025206988720857
$arr[18,9] = 1.033567183548713
$arr[18,10] = 1.02360206886222
$arr[18,11] = 1.044052521841604
$arr[18,12] = 1.012464131705078
$arr[19,0] = 1.02
$arr[19,1] = 1.017968793506468
$arr[19,2] = 1.029481891430291
$arr[19,3] = 1.019394888861164
$arr[19,4] = 1.039822868270547
$arr[19,5] = 1
$arr[19,7] = 1.031552094910182
$arr[19,8] = 1.024687123122871
$arr[19,9] = 1.033086457472029
$arr[19,10] = 1.023037940939669
$arr[19,11] = 1.043388837055312
$arr[19,12] = 1.012289993548533
$arr[20,0] = 1.02
$arr[20,1] = 1.017388301117917
$arr[20,2] = 1.029047843530632
$arr[20,3] = 1.018907290791854
$arr[20,4] = 1.039275282108509
$arr[20,5] = 1
$arr[20,7] = 1.031466059632772
$arr[20,8] = 1.024360434463222
$arr[20,9] = 1.032784032557225
$arr[20,10] = 1.022683701817818
$arr[20,11] = 1.042971786992706
$arr[20,12] = 1.012180517175136
$arr[21,0] = 1.02
$arr[21,1] = 1.017695948035722
$arr[21,2] = 1.029277884879417
$arr[21,3] = 1.019165662917179
$arr[21,4] = 1.039565477282077
$arr[21,5] = 1
$arr[21,7] = 1.031511780115581
$arr[21,8] = 1.024533601594426
$arr[21,9] = 1.032944369674197
$arr[21,10] = 1.022871447482665
$arr[21,11] = 1.043192850239801
$arr[21,12] = 1.012238551444045
$arr[22,0] = 1.02
$arr[22,1] = 1.018908957494647
$arr[22,2] = 1.030184759030797
$arr[22,3] = 1.020185340840555
$arr[22,4] = 1.040709945269298
$arr[22,5] = 1
$arr[22,7] = 1.031689299372728
$arr[22,8] = 1.025215707142779
$arr[22,9] = 1.033575239939262
$arr[22,10] = 1.023611534068929
$arr[22,11] = 1.044063652450792
$arr[22,12] = 1.012467051307268
$arr[23,0] = 1.02
$arr[23,1] = 1.020320316900033
$arr[23,2] = 1.031239595551549
$arr[23,3] = 1.021373707075015
$arr[23,4] = 1.042042090272557
$arr[23,5] = 1
$arr[23,7] = 1.031890141023106
$arr[23,8] = 1.026007953567774
$arr[23,9] = 1.034306526200862
$arr[23,10] = 1.024472275158266
$arr[23,11] = 1.045075130797477
$arr[23,12] = 1.012732244652742

$ws.Range("B2:N25").Value = $arr